# Translate the "is_positive" column values from Turkish to English
# across every worksheet in the workbook:
#   "Hayır" -> "No"
#   "Evet"  -> "Yes"
#
# The header cell in row 1 (literal text "is_positive") does not contain
# either Turkish word, so it is left untouched by these replacements.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    [void]$usedRange.Replace("Hayır", "No")
    [void]$usedRange.Replace("Evet", "Yes")
}
